# Separate xamls for login, transaction and reading the data
# Populate the "Settings" sheet (first sheet) with the new OutputDataFilePath /
# column-mapping configuration rows, and move the active selection to B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value = "OutputDataFilePath"
$ws.Range("B8").Value = "Sheet2"

$ws.Range("B9").Value = "Status"
$ws.Range("B10").Value = "Exception"
$ws.Range("B11").Value = "TransactionNumber"

$ws.Range("A9").Value = "StatusColumn"
$ws.Range("A10").Value = "ExceptionColumn"
$ws.Range("A11").Value = "TransactionNumberColumn"

$ws.Range("B8").Select()
